$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 22.79715566666667
$ws.Range("H2").Value = 68.391467
$ws.Range("I2").Value = 0.181315435549823
$ws.Range("J2").Value = 0.181315435549823
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.490547666666667
$ws.Range("N2").Value = 19.471643
$ws.Range("O2").Value = 0.8021666724616637
$ws.Range("P2").Value = 0.8021666724616636
$ws.Range("Q2").Value = 147.9660255189201
$ws.Range("R2").Value = 1331.694229670281
$ws.Range("S2").Value = 0.1454451996009388
$ws.Range("T2").Value = 0.1454451996009388

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 22.79715566666667
$ws.Range("H3").Value = 68.391467
$ws.Range("I3").Value = 0.181315435549823
$ws.Range("J3").Value = 0.181315435549823
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.600723
$ws.Range("N3").Value = 4.802169
$ws.Range("O3").Value = 0.1978333275383364
$ws.Range("P3").Value = 0.1978333275383364
$ws.Range("Q3").Value = 36.49193141021367
$ws.Range("R3").Value = 328.427382691923
$ws.Range("S3").Value = 0.03587023594888426
$ws.Range("T3").Value = 0.03587023594888426

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 51.67462766666667
$ws.Range("H4").Value = 155.023883
$ws.Range("I4").Value = 0.4109902024293441
$ws.Range("J4").Value = 0.4109902024293441
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.490547666666667
$ws.Range("N4").Value = 19.471643
$ws.Range("O4").Value = 0.8021666724616637
$ws.Range("P4").Value = 0.8021666724616636
$ws.Range("Q4").Value = 335.3966340277522
$ws.Range("R4").Value = 3018.569706249769
$ws.Range("S4").Value = 0.3296826430970926
$ws.Range("T4").Value = 0.3296826430970925

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 51.67462766666667
$ws.Range("H5").Value = 155.023883
$ws.Range("I5").Value = 0.4109902024293441
$ws.Range("J5").Value = 0.4109902024293441
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.600723
$ws.Range("N5").Value = 4.802169
$ws.Range("O5").Value = 0.1978333275383364
$ws.Range("P5").Value = 0.1978333275383364
$ws.Range("Q5").Value = 82.71676502246969
$ws.Range("R5").Value = 744.450885202227
$ws.Range("S5").Value = 0.08130755933225162
$ws.Range("T5").Value = 0.0813075593322516

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 23.76107833333333
$ws.Range("H6").Value = 71.28323499999999
$ws.Range("I6").Value = 0.1889819208209905
$ws.Range("J6").Value = 0.1889819208209905
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 6.490547666666667
$ws.Range("N6").Value = 19.471643
$ws.Range("O6").Value = 0.8021666724616637
$ws.Range("P6").Value = 0.8021666724616636
$ws.Range("Q6").Value = 154.2224115339006
$ws.Range("R6").Value = 1388.001703805105
$ws.Range("S6").Value = 0.1515949985803876
$ws.Range("T6").Value = 0.1515949985803876

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 23.76107833333333
$ws.Range("H7").Value = 71.28323499999999
$ws.Range("I7").Value = 0.1889819208209905
$ws.Range("J7").Value = 0.1889819208209905
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.600723
$ws.Range("N7").Value = 4.802169
$ws.Range("O7").Value = 0.1978333275383364
$ws.Range("P7").Value = 0.1978333275383364
$ws.Range("Q7").Value = 38.03490459296833
$ws.Range("R7").Value = 342.314141336715
$ws.Range("S7").Value = 0.03738692224060298
$ws.Range("T7").Value = 0.03738692224060297

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 27.49915666666666
$ws.Range("H8").Value = 82.49746999999999
$ws.Range("I8").Value = 0.2187124411998423
$ws.Range("J8").Value = 0.2187124411998423
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.490547666666667
$ws.Range("N8").Value = 19.471643
$ws.Range("O8").Value = 0.8021666724616637
$ws.Range("P8").Value = 0.8021666724616636
$ws.Range("Q8").Value = 178.4845871381344
$ws.Range("R8").Value = 1606.36128424321
$ws.Range("S8").Value = 0.1754438311832448
$ws.Range("T8").Value = 0.1754438311832448

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 27.49915666666666
$ws.Range("H9").Value = 82.49746999999999
$ws.Range("I9").Value = 0.2187124411998423
$ws.Range("J9").Value = 0.2187124411998423
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.600723
$ws.Range("N9").Value = 4.802169
$ws.Range("O9").Value = 0.1978333275383364
$ws.Range("P9").Value = 0.1978333275383364
$ws.Range("Q9").Value = 44.01853255693666
$ws.Range("R9").Value = 396.16679301243
$ws.Range("S9").Value = 0.04326861001659755
$ws.Range("T9").Value = 0.04326861001659754

